$d = $word.ActiveDocument

# Locate the paragraph that ends with the "You might want to replace
# m: endif by m:endif" validation message (the m:endif diagnostic
# paragraph). We search by content rather than a hard-coded index so the
# script is resilient to minor structural differences.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*You might want to replace m: endif by m:endif*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    # Append three new runs at the very end of the paragraph, right
    # before the paragraph mark:
    #   1) "    "                (plain run, 4 spaces)
    #   2) "<---"                (red, 16pt, light-gray highlight)
    #   3) "missing expression"  (red, 16pt, light-gray highlight)

    # 1) plain spacer run
    $pEnd = $target.Range.End
    $insertPoint = $d.Range($pEnd - 1, $pEnd - 1)
    $insertPoint.InsertAfter("    ")

    # 2) "<---" marker run
    $pEnd = $target.Range.End
    $insertStart = $pEnd - 1
    $insertPoint = $d.Range($insertStart, $insertStart)
    $insertPoint.InsertAfter("<---")
    $markerRange = $d.Range($insertStart, $insertStart + 4)
    $markerRange.Font.Color = 255
    $markerRange.Font.Size = 16
    $markerRange.Font.HighlightColorIndex = 16

    # 3) "missing expression" message run
    $pEnd = $target.Range.End
    $insertStart = $pEnd - 1
    $insertPoint = $d.Range($insertStart, $insertStart)
    $insertPoint.InsertAfter("missing expression")
    $msgRange = $d.Range($insertStart, $insertStart + 19)
    $msgRange.Font.Color = 255
    $msgRange.Font.Size = 16
    $msgRange.Font.HighlightColorIndex = 16
}
